# Update data and data checking rules
# Adds a new date row (2020-04-14) to both the "Confirmados" and "Mortes"
# sheets, mirroring the existing row layout (row 51, columns A:AB).

function Add-DailyRow($SheetName, $Row, $DateText, $Values) {
    $ws = $wb.Worksheets.Item($SheetName)

    # Write the date as literal text (matching the existing column A cells,
    # which are stored as shared strings, not date serials). Entering the
    # text directly would be auto-recognised as a date by Excel, so we
    # round-trip it through a TEXT() formula and paste the computed result
    # back as a plain value/text cell.
    $dateCell = $ws.Range("A$Row")
    $dateCell.Formula = '=TEXT("' + $DateText + '","@")'
    $dateCell.Copy()
    $dateCell.PasteSpecial(-4163)
    $excel.CutCopyMode = $false

    $col = 2
    foreach ($val in $Values) {
        $ws.Cells.Item($Row, $col).Value = $val
        $col++
    }
}

$wb = $excel.ActiveWorkbook

# Values for the new row, in column order B..AB (states in the same order
# as the existing header row: Acre, Alagoas, Amapa, Amazonas, Bahia, Ceara,
# Distrito Federal, Espirito Santo, Goias, Maranhao, Mato Grosso,
# Mato Grosso do Sul, Minas Gerais, Para, Paraiba, Parana, Pernambuco,
# Piaui, Rio de Janeiro, Rio Grande do Norte, Rio Grande do Sul, Rondonia,
# Roraima, Santa Catarina, Sao Paulo, Sergipe, Tocantins)
$confirmados = @(99,72,307,1484,759,2005,651,463,284,478,138,115,884,323,136,791,1284,58,3410,376,700,64,113,826,9371,45,26)
$mortes      = @(3,4,6,90,22,107,17,17,15,32,4,4,27,19,16,36,115,8,224,18,18,2,3,26,695,4,0)

$newDate = "2020-04-14"
$newRow = 51

Add-DailyRow "Confirmados" $newRow $newDate $confirmados
Add-DailyRow "Mortes" $newRow $newDate $mortes

Write-Host "Added $newDate row to Confirmados and Mortes sheets"
